$d = $word.ActiveDocument

# The checklist header table (Project Name / Sprint No. / Review Date / File Name)
# is the first table in the document.
$t = $d.Tables.Item(1)

# "Sprint No." value cell: row 2, column 4 -> change "1" to "2"
$sprintCell = $t.Cell(2, 4)
$sprintCell.Range.Text = "2"

# "Review Date" value cell: row 3, column 2 -> change "02/09/18" to "02/21/18"
$dateCell = $t.Cell(3, 2)
$dateCell.Range.Text = "02/21/18"
